# PCB and base changes
# Updates the Arduino Nano BOM: splits header-pin line items into explicit
# male/female 1x3/1x4/1x15 header rows, adds a Notes column header, clarifies
# the bottle-cap note, and relocates the ultrasonic sensor note into the
# Notes column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new "1x15 Female header" line item -------------------
# Inserting a whole row at 14 shifts "Continuous Servo" (old row 14) and
# everything below it down by one, automatically re-pointing the shared
# D-column formulas and the trailing SUM() total.
$ws.Rows.Item(14).Insert()

# --- New "Notes" column header ----------------------------------------------
# Match the formatting of the rest of the header row (D2) before typing text.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E2").Value = "Notes"

# --- Row 3: bottle cap note clarified ---------------------------------------
$ws.Range("E3").Value = "Buy your own, Gatorade bottles work best"

# --- Row 12: "Male header Pins" (6 @ $0.02) -> "1x3 Male header " (2 @ $0.06)
$ws.Range("A12").Value = "1x3 Male header "
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 0.06

# --- Row 13: "Female Header Pins" (80 @ $0.02) -> "1x4 Female header " (1 @ $0.08)
$ws.Range("A13").Value = "1x4 Female header "
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0.08

# --- Row 14 (new): "1x15 Female header " (2 @ $0.30) with a note -----------
$ws.Range("A14").Value = "1x15 Female header "
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 0.3
$ws.Range("D14").Formula = "=B14*C14"
$ws.Range("E14").Value = "You might need to buy a larger one and then break off the pins you need"

# --- Relocate the ultrasonic-sensor note from column F to column E ---------
# It travelled with "Sensor Sonar" during the row insert (old row 16 -> 17);
# move it from F17 into the new Notes column E17.
$ws.Range("E17").Value = $ws.Range("F17").Value2
$ws.Range("F17").ClearContents()

# --- Refresh the sheet's selection to cover the (now one-row-taller) table -
$ws.Range("A1:E22").Select() | Out-Null

$excel.Calculate()
